# Add two new columns "I0" (col I) and "IF" (col J) to the header row,
# matching the style of the existing header cells, and fill in their
# values for the single data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (bold font, borders, centered alignment) from the
# last existing header cell (H1) onto the two new header cells so they
# look consistent with the rest of the header row.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Header row (row 1): new column labels.
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data row (row 2): new values, left unstyled like the other data cells.
$ws.Range("I2").Value = 3
$ws.Range("J2").Value = 7
